# The first paragraph of the room-printout template ("Gebäude #g / Raum #r")
# mixes two font sizes: the paragraph mark (and a hidden noProof drawing run)
# is 30pt (sz=60) while the visible text runs are 36pt (sz=72). The layout
# fix unifies everything in that paragraph to 32pt (sz/szCs = 64).
#
# Font.Size drives w:sz (ascii/east-asian size) and Font.SizeBi drives
# w:szCs (complex-script size); Word keeps both in sync when you touch a
# whole paragraph range, so setting both on the paragraph's Range updates
# every run (including the paragraph-mark's own rPr) in one shot.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range

$r.Font.Size = 32
$r.Font.SizeBi = 32
